$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 3.204747333333333
$ws.Range("H2").Value = 9.614241999999999
$ws.Range("I2").Value = 0.01973032100547387
$ws.Range("J2").Value = 0.01973032100547387
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 31.82741333333333
$ws.Range("N2").Value = 95.48223999999999
$ws.Range("O2").Value = 0.114390792932228
$ws.Range("P2").Value = 0.114390792932228
$ws.Range("Q2").Value = 101.9988180068978
$ws.Range("R2").Value = 917.9893620620799
$ws.Range("S2").Value = 0.00225696706462355
$ws.Range("T2").Value = 0.002256967064623551
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 3.204747333333333
$ws.Range("H3").Value = 9.614241999999999
$ws.Range("I3").Value = 0.01973032100547387
$ws.Range("J3").Value = 0.01973032100547387
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 85.46317833333335
$ws.Range("N3").Value = 256.389535
$ws.Range("O3").Value = 0.307162904935779
$ws.Range("P3").Value = 0.307162904935779
$ws.Range("Q3").Value = 273.8878928619411
$ws.Range("R3").Value = 2464.99103575747
$ws.Range("S3").Value = 0.006060422715356776
$ws.Range("T3").Value = 0.006060422715356776
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 3.204747333333333
$ws.Range("H4").Value = 9.614241999999999
$ws.Range("I4").Value = 0.01973032100547387
$ws.Range("J4").Value = 0.01973032100547387
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 122.2478306666667
$ws.Range("N4").Value = 366.743492
$ws.Range("O4").Value = 0.4393704929064738
$ws.Range("P4").Value = 0.4393704929064738
$ws.Range("Q4").Value = 391.7734093347849
$ws.Range("R4").Value = 3525.960684013064
$ws.Range("S4").Value = 0.008668920865378009
$ws.Range("T4").Value = 0.008668920865378011
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 3.204747333333333
$ws.Range("H5").Value = 9.614241999999999
$ws.Range("I5").Value = 0.01973032100547387
$ws.Range("J5").Value = 0.01973032100547387
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 38.69562533333333
$ws.Range("N5").Value = 116.086876
$ws.Range("O5").Value = 0.1390758092255191
$ws.Range("P5").Value = 0.1390758092255191
$ws.Range("Q5").Value = 124.0097020986658
$ws.Range("R5").Value = 1116.087318887992
$ws.Range("S5").Value = 0.002744010360115537
$ws.Range("T5").Value = 0.002744010360115537
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 145.2141163333334
$ws.Range("H6").Value = 435.6423490000001
$ws.Range("I6").Value = 0.8940240311559332
$ws.Range("J6").Value = 0.8940240311559333
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 31.82741333333333
$ws.Range("N6").Value = 95.48223999999999
$ws.Range("O6").Value = 0.114390792932228
$ws.Range("P6").Value = 0.114390792932228
$ws.Range("Q6").Value = 4621.789702375751
$ws.Range("R6").Value = 41596.10732138177
$ws.Range("S6").Value = 0.1022681178243941
$ws.Range("T6").Value = 0.1022681178243942
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 145.2141163333334
$ws.Range("H7").Value = 435.6423490000001
$ws.Range("I7").Value = 0.8940240311559332
$ws.Range("J7").Value = 0.8940240311559333
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 85.46317833333335
$ws.Range("N7").Value = 256.389535
$ws.Range("O7").Value = 0.307162904935779
$ws.Range("P7").Value = 0.307162904935779
$ws.Range("Q7").Value = 12410.45992071308
$ws.Range("R7").Value = 111694.1392864177
$ws.Range("S7").Value = 0.2746110184922518
$ws.Range("T7").Value = 0.2746110184922519
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 145.2141163333334
$ws.Range("H8").Value = 435.6423490000001
$ws.Range("I8").Value = 0.8940240311559332
$ws.Range("J8").Value = 0.8940240311559333
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 122.2478306666667
$ws.Range("N8").Value = 366.743492
$ws.Range("O8").Value = 0.4393704929064738
$ws.Range("P8").Value = 0.4393704929064738
$ws.Range("Q8").Value = 17752.11070392697
$ws.Range("R8").Value = 159768.9963353427
$ws.Range("S8").Value = 0.3928077792392151
$ws.Range("T8").Value = 0.3928077792392152
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 145.2141163333334
$ws.Range("H9").Value = 435.6423490000001
$ws.Range("I9").Value = 0.8940240311559332
$ws.Range("J9").Value = 0.8940240311559333
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 38.69562533333333
$ws.Range("N9").Value = 116.086876
$ws.Range("O9").Value = 0.1390758092255191
$ws.Range("P9").Value = 0.1390758092255191
$ws.Range("Q9").Value = 5619.151038745747
$ws.Range("R9").Value = 50572.35934871173
$ws.Range("S9").Value = 0.1243371156000721
$ws.Range("T9").Value = 0.1243371156000721
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.8052786666666667
$ws.Range("H10").Value = 2.415836
$ws.Range("I10").Value = 0.004957771998726471
$ws.Range("J10").Value = 0.004957771998726472
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 31.82741333333333
$ws.Range("N10").Value = 95.48223999999999
$ws.Range("O10").Value = 0.114390792932228
$ws.Range("P10").Value = 0.114390792932228
$ws.Range("Q10").Value = 25.62993697251555
$ws.Range("R10").Value = 230.66943275264
$ws.Range("S10").Value = 0.000567123470111518
$ws.Range("T10").Value = 0.0005671234701115182
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.8052786666666667
$ws.Range("H11").Value = 2.415836
$ws.Range("I11").Value = 0.004957771998726471
$ws.Range("J11").Value = 0.004957771998726472
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 85.46317833333335
$ws.Range("N11").Value = 256.389535
$ws.Range("O11").Value = 0.307162904935779
$ws.Range("P11").Value = 0.307162904935779
$ws.Range("Q11").Value = 68.82167429736224
$ws.Range("R11").Value = 619.3950686762601
$ws.Range("S11").Value = 0.001522843649138086
$ws.Range("T11").Value = 0.001522843649138086
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.8052786666666667
$ws.Range("H12").Value = 2.415836
$ws.Range("I12").Value = 0.004957771998726471
$ws.Range("J12").Value = 0.004957771998726472
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 122.2478306666667
$ws.Range("N12").Value = 366.743492
$ws.Range("O12").Value = 0.4393704929064738
$ws.Range("P12").Value = 0.4393704929064738
$ws.Range("Q12").Value = 98.44357008214578
$ws.Range("R12").Value = 885.9921307393121
$ws.Range("S12").Value = 0.002178298726798363
$ws.Range("T12").Value = 0.002178298726798364
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.8052786666666667
$ws.Range("H13").Value = 2.415836
$ws.Range("I13").Value = 0.004957771998726471
$ws.Range("J13").Value = 0.004957771998726472
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 38.69562533333333
$ws.Range("N13").Value = 116.086876
$ws.Range("O13").Value = 0.1390758092255191
$ws.Range("P13").Value = 0.1390758092255191
$ws.Range("Q13").Value = 31.16076157425956
$ws.Range("R13").Value = 280.446854168336
$ws.Range("S13").Value = 0.0006895061526785033
$ws.Range("T13").Value = 0.0006895061526785034
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 13.203389
$ws.Range("H14").Value = 39.610167
$ws.Range("I14").Value = 0.08128787583986632
$ws.Range("J14").Value = 0.08128787583986634
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 31.82741333333333
$ws.Range("N14").Value = 95.48223999999999
$ws.Range("O14").Value = 0.114390792932228
$ws.Range("P14").Value = 0.114390792932228
$ws.Range("Q14").Value = 420.2297191037866
$ws.Range("R14").Value = 3782.06747193408
$ws.Range("S14").Value = 0.009298584573098808
$ws.Range("T14").Value = 0.009298584573098812
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 13.203389
$ws.Range("H15").Value = 39.610167
$ws.Range("I15").Value = 0.08128787583986632
$ws.Range("J15").Value = 0.08128787583986634
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 85.46317833333335
$ws.Range("N15").Value = 256.389535
$ws.Range("O15").Value = 0.307162904935779
$ws.Range("P15").Value = 0.307162904935779
$ws.Range("Q15").Value = 1128.403588711372
$ws.Range("R15").Value = 10155.63229840235
$ws.Range("S15").Value = 0.02496862007903227
$ws.Range("T15").Value = 0.02496862007903227
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 13.203389
$ws.Range("H16").Value = 39.610167
$ws.Range("I16").Value = 0.08128787583986632
$ws.Range("J16").Value = 0.08128787583986634
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 122.2478306666667
$ws.Range("N16").Value = 366.743492
$ws.Range("O16").Value = 0.4393704929064738
$ws.Range("P16").Value = 0.4393704929064738
$ws.Range("Q16").Value = 1614.08566269813
$ws.Range("R16").Value = 14526.77096428317
$ws.Range("S16").Value = 0.03571549407508231
$ws.Range("T16").Value = 0.03571549407508232
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 13.203389
$ws.Range("H17").Value = 39.610167
$ws.Range("I17").Value = 0.08128787583986632
$ws.Range("J17").Value = 0.08128787583986634
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 38.69562533333333
$ws.Range("N17").Value = 116.086876
$ws.Range("O17").Value = 0.1390758092255191
$ws.Range("P17").Value = 0.1390758092255191
$ws.Range("Q17").Value = 510.9133938742547
$ws.Range("R17").Value = 4598.220544868292
$ws.Range("S17").Value = 0.01130517711265294
$ws.Range("T17").Value = 0.01130517711265294
